# Rotate the comma-separated "Recorded By" names in column G so that the
# first name in the list is moved to the end of the list.
# e.g. "backup@backdoor.com, System" -> "System, backup@backdoor.com"
# Cells that only contain a single name are left untouched (rotation is a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    if ($parts.Count -gt 1) {
        $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
        $cell.Value2 = $rotated
    }
}
